# Append 2021 and 2022 rows (rows 11 and 12) to the "资产总计" (total assets)
# time series sheet, matching the style already used by the existing year
# label cells in column A (copy format from A10, then overwrite the value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 50541.91
$ws.Range("C11").Value = 13173.28
$ws.Range("D11").Value = 3542.95
$ws.Range("E11").Value = 8.57
$ws.Range("F11").Value = 33421.95
$ws.Range("G11").Value = 87709.98
$ws.Range("H11").Value = 10278.42
$ws.Range("I11").Value = 44053.13
$ws.Range("J11").Value = 7012.61
$ws.Range("K11").Value = 6957.59
$ws.Range("L11").Value = 4666.65
$ws.Range("M11").Value = 3123.5
$ws.Range("N11").Value = 9561.610000000001
$ws.Range("O11").Value = 44470.72
$ws.Range("P11").Value = 6598
$ws.Range("Q11").Value = 5689.64
$ws.Range("R11").Value = 27449.42
$ws.Range("S11").Value = 24077.29
$ws.Range("T11").Value = 89241.28
$ws.Range("U11").Value = 11702.64
$ws.Range("V11").Value = 68700.96000000001
$ws.Range("W11").Value = 15439.26
$ws.Range("X11").Value = 202577.01
$ws.Range("Y11").Value = 89975.35000000001
$ws.Range("Z11").Value = 7045.57
$ws.Range("AA11").Value = 41344.07
$ws.Range("AB11").Value = 22840.76
$ws.Range("AC11").Value = 21840.83
$ws.Range("AD11").Value = 11633.91
$ws.Range("AE11").Value = 1466716.3
$ws.Range("AF11").Value = 153315.61
$ws.Range("AG11").Value = 56217.77
$ws.Range("AH11").Value = 15922.07
$ws.Range("AI11").Value = 20521.25
$ws.Range("AJ11").Value = 2603.99
$ws.Range("AK11").Value = 40861.98
$ws.Range("AL11").Value = 31477.42
$ws.Range("AM11").Value = 73304.25
$ws.Range("AN11").Value = 7146.36
$ws.Range("AO11").Value = 19040.19
$ws.Range("AP11").Value = 69199.82000000001
$ws.Range("AQ11").Value = 12426.72

$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "2022年"
$ws.Range("B12").Value = 55794.4
$ws.Range("C12").Value = 14141.4
$ws.Range("D12").Value = 1724.5
$ws.Range("E12").Value = 7.8
$ws.Range("F12").Value = 35943.9
$ws.Range("G12").Value = 95670.10000000001
$ws.Range("H12").Value = 11231
$ws.Range("I12").Value = 47885.3
$ws.Range("J12").Value = 7371
$ws.Range("K12").Value = 7147.3
$ws.Range("L12").Value = 5671.9
$ws.Range("M12").Value = 3248.2
$ws.Range("N12").Value = 9943.200000000001
$ws.Range("O12").Value = 48072.2
$ws.Range("P12").Value = 6738
$ws.Range("Q12").Value = 6039.1
$ws.Range("R12").Value = 28922.4
$ws.Range("S12").Value = 26482.2
$ws.Range("T12").Value = 97297.8
$ws.Range("U12").Value = 10658.7
$ws.Range("V12").Value = 73210.5
$ws.Range("W12").Value = 16435.4
$ws.Range("X12").Value = 213137
$ws.Range("Y12").Value = 108288.2
$ws.Range("Z12").Value = 6984.4
$ws.Range("AA12").Value = 44247.8
$ws.Range("AB12").Value = 24114.8
$ws.Range("AC12").Value = 22508.4
$ws.Range("AD12").Value = 11418.6
$ws.Range("AE12").Value = 1561196.7
$ws.Range("AF12").Value = 169812.3
$ws.Range("AG12").Value = 58111.5
$ws.Range("AH12").Value = 16470.1
$ws.Range("AI12").Value = 22641.5
$ws.Range("AJ12").Value = 2692.5
$ws.Range("AK12").Value = 38804.3
$ws.Range("AL12").Value = 19767.1
$ws.Range("AM12").Value = 78880.3
$ws.Range("AN12").Value = 8042.6
$ws.Range("AO12").Value = 21304.3
$ws.Range("AP12").Value = 71800.3
$ws.Range("AQ12").Value = 12534.2
